$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2533.3333
$ws.Range("I40").Value = 2400.1428
$ws.Range("J40").Value = 2999.5
$ws.Range("K40").Value = 2400.1428
$ws.Range("L40").Value = 2999.5
$ws.Range("M40").Value = -2225.1428
$ws.Range("N40").Value = -3349.5
$ws.Range("H76").Value = 3634.3333
$ws.Range("I76").Value = 3951.5
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 3951.5
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -3636.5
$ws.Range("N76").Value = -3630
$ws.Range("H79").Value = 3634.3333
$ws.Range("I79").Value = 3951.5
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 3951.5
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -2859.5
$ws.Range("N79").Value = -5184
$ws.Range("H132").Value = 9011756
$ws.Range("I132").Value = 12348577
$ws.Range("K132").Value = 37045731
$ws.Range("M132").Value = -37043201
$ws.Range("H133").Value = 35796
$ws.Range("J133").Value = 35796
$ws.Range("L133").Value = 35796
$ws.Range("N133").Value = -45916
$ws.Range("H135").Value = 151.52632
$ws.Range("I135").Value = 114.9375
$ws.Range("K135").Value = 1034.4375
$ws.Range("M135").Value = 1500.5625
$ws.Range("H138").Value = 1821.22
$ws.Range("I138").Value = 1057.8334
$ws.Range("J138").Value = 1869.9468
$ws.Range("K138").Value = 3173.5002
$ws.Range("L138").Value = 5609.8404
$ws.Range("M138").Value = 1966.4998
$ws.Range("N138").Value = -15889.8404

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2822.3728
$ws.Range("I32").Value = 3203.318
$ws.Range("J32").Value = 1704.9333
$ws.Range("K32").Value = 3203.318
$ws.Range("L32").Value = 1704.9333
$ws.Range("M32").Value = -2916.318
$ws.Range("N32").Value = -2278.9333
$ws.Range("H61").Value = 1420.2
$ws.Range("I61").Value = 1420.2
$ws.Range("K61").Value = 1420.2
$ws.Range("M61").Value = -1208.2
$ws.Range("H74").Value = 1136.3334
$ws.Range("I74").Value = 983.4375
$ws.Range("K74").Value = 983.4375
$ws.Range("M74").Value = -109.4375
$ws.Range("H77").Value = 1136.3334
$ws.Range("I77").Value = 983.4375
$ws.Range("K77").Value = 4917.1875
$ws.Range("M77").Value = -549.1875
$ws.Range("H122").Value = 1399.7059
$ws.Range("I122").Value = 1319.7333
$ws.Range("K122").Value = 3959.199900000001
$ws.Range("M122").Value = -1509.199900000001
$ws.Range("H123").Value = 66685.8
$ws.Range("J123").Value = 66685.8
$ws.Range("L123").Value = 66685.8
$ws.Range("N123").Value = -76485.8
$ws.Range("H132").Value = 3155.6428
$ws.Range("I132").Value = 2928.087
$ws.Range("J132").Value = 4202.4
$ws.Range("K132").Value = 8784.261
$ws.Range("L132").Value = 12607.2
$ws.Range("M132").Value = -6254.261
$ws.Range("N132").Value = -17667.2
$ws.Range("H133").Value = 25717.143
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 25717.143
$ws.Range("K133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("M133").Value = 25717.143
$ws.Range("N133").Value = -30777.143
$ws.Range("H136").Value = 1420.2
$ws.Range("I136").Value = 1420.2
$ws.Range("K136").Value = 4260.6
$ws.Range("M136").Value = -1710.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 11663.81
$ws.Range("I134").Value = 7247.05
$ws.Range("K134").Value = 21741.15
$ws.Range("M134").Value = -19206.15

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1408
$ws.Range("I58").Value = 1524.5714
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 1524.5714
$ws.Range("L58").Value = 1000
$ws.Range("M58").Value = -1321.5714
$ws.Range("N58").Value = -1406
$ws.Range("H99").Value = 1463551.2
$ws.Range("I99").Value = 2633032.8
$ws.Range("K99").Value = 2633032.8
$ws.Range("M99").Value = -2631534.8
$ws.Range("H126").Value = 1463551.2
$ws.Range("I126").Value = 2633032.8
$ws.Range("K126").Value = 7899098.399999999
$ws.Range("M126").Value = -7896628.399999999
$ws.Range("H136").Value = 1408
$ws.Range("I136").Value = 1524.5714
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 4573.7142
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -2023.7142
$ws.Range("N136").Value = -8100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1608.3024
$ws.Range("J68").Value = 1860.4412
$ws.Range("L68").Value = 5581.3236
$ws.Range("N68").Value = -7203.3236
$ws.Range("H71").Value = 1608.3024
$ws.Range("J71").Value = 1860.4412
$ws.Range("L71").Value = 16743.9708
$ws.Range("N71").Value = -24855.9708
$ws.Range("H140").Value = 29526.055
$ws.Range("J140").Value = 2869.1667
$ws.Range("L140").Value = 8607.500100000001
$ws.Range("N140").Value = -18967.5001
$ws.Range("H141").Value = 1901.3334
$ws.Range("I141").Value = 1901.3334
$ws.Range("K141").Value = 5704.0002
$ws.Range("M141").Value = -524.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3877.6667
$ws.Range("I80").Value = 2589.8
$ws.Range("J80").Value = 5487.5
$ws.Range("K80").Value = 2589.8
$ws.Range("L80").Value = 5487.5
$ws.Range("M80").Value = -1591.8
$ws.Range("N80").Value = -7483.5
$ws.Range("H83").Value = 3877.6667
$ws.Range("I83").Value = 2589.8
$ws.Range("J83").Value = 5487.5
$ws.Range("K83").Value = 12949
$ws.Range("L83").Value = 27437.5
$ws.Range("M83").Value = -7957
$ws.Range("N83").Value = -37421.5
$ws.Range("H126").Value = 2966.1365
$ws.Range("I126").Value = 1809.091
$ws.Range("J126").Value = 4123.1816
$ws.Range("K126").Value = 5427.272999999999
$ws.Range("L126").Value = 12369.5448
$ws.Range("M126").Value = -2957.272999999999
$ws.Range("N126").Value = -17309.5448
$ws.Range("H132").Value = 2917.4736
$ws.Range("I132").Value = 2648.7693
$ws.Range("J132").Value = 3499.6667
$ws.Range("K132").Value = 7946.3079
$ws.Range("L132").Value = 10499.0001
$ws.Range("M132").Value = -5416.3079
$ws.Range("N132").Value = -15559.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1312.4375
$ws.Range("I16").Value = 1088.3846
$ws.Range("J16").Value = 2283.3333
$ws.Range("K16").Value = 1088.3846
$ws.Range("L16").Value = 2283.3333
$ws.Range("M16").Value = -918.3846000000001
$ws.Range("N16").Value = -2623.3333
$ws.Range("H68").Value = 1481.909
$ws.Range("I68").Value = 1171.8572
$ws.Range("J68").Value = 2024.5
$ws.Range("K68").Value = 1171.8572
$ws.Range("L68").Value = 2024.5
$ws.Range("M68").Value = -422.8571999999999
$ws.Range("N68").Value = -3522.5
$ws.Range("H71").Value = 1481.909
$ws.Range("I71").Value = 1171.8572
$ws.Range("J71").Value = 2024.5
$ws.Range("K71").Value = 5859.286
$ws.Range("L71").Value = 10122.5
$ws.Range("M71").Value = -2115.286
$ws.Range("N71").Value = -17610.5
$ws.Range("H82").Value = 2972.6365
$ws.Range("I82").Value = 3000
$ws.Range("J82").Value = 2924.75
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 2924.75
$ws.Range("M82").Value = -2639
$ws.Range("N82").Value = -3646.75
$ws.Range("H85").Value = 2972.6365
$ws.Range("I85").Value = 3000
$ws.Range("J85").Value = 2924.75
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 2924.75
$ws.Range("M85").Value = -1752
$ws.Range("N85").Value = -5420.75
$ws.Range("H100").Value = 1975.8334
$ws.Range("I100").Value = 1972.1428
$ws.Range("K100").Value = 1972.1428
$ws.Range("M100").Value = -1431.1428
$ws.Range("H122").Value = 25761666
$ws.Range("I122").Value = 35420416
$ws.Range("K122").Value = 106261248
$ws.Range("M122").Value = -106258798
$ws.Range("H132").Value = 103489.1
$ws.Range("I132").Value = 4629.6665
$ws.Range("J132").Value = 145857.42
$ws.Range("K132").Value = 13888.9995
$ws.Range("L132").Value = 437572.26
$ws.Range("M132").Value = -11358.9995
$ws.Range("N132").Value = -442632.26
$ws.Range("H136").Value = 6678.4443
$ws.Range("I136").Value = 12133.111
$ws.Range("J136").Value = 1223.7778
$ws.Range("K136").Value = 36399.333
$ws.Range("L136").Value = 3671.3334
$ws.Range("M136").Value = -33849.333
$ws.Range("N136").Value = -8771.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 20000
$ws.Range("J54").Value = 20000
$ws.Range("L54").Value = 20000
$ws.Range("N54").Value = -21040
$ws.Range("H123").Value = 56252.332
$ws.Range("J123").Value = 56252.332
$ws.Range("L123").Value = 56252.332
$ws.Range("N123").Value = -66052.33199999999
$ws.Range("H125").Value = 68298.75
$ws.Range("J125").Value = 68298.75
$ws.Range("L125").Value = 68298.75
$ws.Range("N125").Value = -78138.75
$ws.Range("H132").Value = 5571.7144
$ws.Range("I132").Value = 6749.75
$ws.Range("K132").Value = 20249.25
$ws.Range("M132").Value = -17719.25
